# Documentation and work diary - PDF export
# Applies the authoring changes recorded in the commit:
#  - minimize the workbook window
#  - switch the sheet to Page Layout view and update the selection
#  - extend the work-journal table with the 22-Jun (44684) entries
#  - add the "Type" column (B) value "Analyse" to every data row
#  - add a new hyperlink on E12, tweak row heights
#  - set print scale + header/footer for the PDF export

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window: minimize it -------------------------------------
$win = $wb.Windows.Item(1)
$win.WindowState = -4140   # xlMinimized

# --- Sheet view: Page Layout (for print/PDF export) --------------------
$ws.Activate()
$excel.ActiveWindow.View = 3   # xlPageLayoutView

# --- Column B ("Type") = "Analyse" for every existing data row ---------
$ws.Range("B2").Value = "Analyse"
$ws.Range("B3").Value = "Analyse"
$ws.Range("B4").Value = "Analyse"
$ws.Range("B5").Value = "Analyse"
$ws.Range("B6").Value = "Analyse"
$ws.Range("B7").Value = "Analyse"
$ws.Range("B8").Value = "Analyse"
$ws.Range("B9").Value = "Analyse"
$ws.Range("B10").Value = "Analyse"

# Row 10 gains extra height (matches manual resize in the source file)
$ws.Rows.Item(10).RowHeight = 30

# --- New rows 11-17 (work on 2022-05-03, serial 44684) ------------------
# Row 11
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 44684
$ws.Range("B11").Value = "Analyse"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Remplissage du fichier de documentation"

# Row 12 (with hyperlink to the IceScrum sandbox, like E9)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 44684
$ws.Range("B12").Value = "Analyse"
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = "Remplissage de la sandbox sur Scrum"
$ws.Range("E12").Value = "https://icescrum.cpnv.ch/p/NEWSWEBSIT/"
[void]$ws.Hyperlinks.Add($ws.Range("E12"), "https://icescrum.cpnv.ch/p/NEWSWEBSIT/")
$ws.Range("E9").Copy()
$ws.Range("E12").PasteSpecial(-4122)

# Row 13
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 44684
$ws.Range("B13").Value = "Analyse"
$ws.Range("C13").Value = 1.5
$ws.Range("D13").Value = "Remplissage du fichier de documentation - Cahier des charges et plus"
$ws.Range("F13").Value = "39d1f2fc5b7a74535261ae7d8b8e759f4d303ebc"
$ws.Rows.Item(13).RowHeight = 30

# Row 14
$ws.Range("A9").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 44684
$ws.Range("B14").Value = "Analyse"
$ws.Range("C14").Value = 0.75
$ws.Range("D14").Value = "Revue de la sandbox sur IceScrum avec CdP"

# Row 15
$ws.Range("A9").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 44684
$ws.Range("B15").Value = "Analyse"
$ws.Range("C15").Value = 0.75
$ws.Range("D15").Value = "Documentation de l'analyse concurentielle"
$ws.Range("F15").Value = "39d1f2fc5b7a74535261ae7d8b8e759f4d303ebc"

# Row 16
$ws.Range("A9").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 44684
$ws.Range("B16").Value = "Analyse"
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = "Définition de tests dans les story sur IceScrum"

# Row 17
$ws.Range("A9").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 44684
$ws.Range("B17").Value = "Analyse"
$ws.Range("C17").Value = 0.25
$ws.Range("D17").Value = "Préparation du premier rendu"
$ws.Range("E17").Value = "J'ai oublié de faire le rendu de la planification initiale le premier soir"
$ws.Rows.Item(17).RowHeight = 30

# --- Selection left where the author left it ---------------------------
[void]$ws.Range("E13").Select()

# --- Print setup: scale for the PDF export + header/footer -------------
$ps = $ws.PageSetup
$ps.Zoom = 65
$ps.LeftHeader = "CPNV"
$ps.LeftFooter = "Louis Richard - louis.richard@cpnv.ch"
$ps.CenterFooter = "&P/&N"
$ps.RightFooter = "&D"
